$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Support raffle claiming by date range:
# Collapse the old standalone "TRANSACTION DATE" row (row 4) and the
# "REFERENCE" / "ENTRIES" row (row 6) into a single row so a date range
# can sit alongside reference/entries on one line.

# Remember the TRANSACTION DATE label before we shuffle rows around.
$txDate = $ws.Range("A4").Value2

# Remove row 4 entirely - this shifts the old row 6 (REFERENCE / ENTRIES) up to row 5.
$ws.Rows(4).Delete()

# Row 5 currently holds: A5 = REFERENCE, B5 = ENTRIES.
# Shift that pair one column to the right (B5:C5), preserving formatting.
$ws.Range("B5").Copy($ws.Range("C5"))
$ws.Range("A5").Copy($ws.Range("B5"))

# Put TRANSACTION DATE back into A5, bold to match the other header cells.
$ws.Range("A5").Value2 = $txDate
$ws.Range("A5").Font.Bold = $true

# New column widths to fit REFERENCE / ENTRIES in their new spots.
# (ColumnWidth is quantized internally to a 1/6-character pixel grid with a
# 5/6 offset, so back the target off by that offset to land on the closest
# representable width to the authored value.)
$gridOffset = 0.8333333333333334
$ws.Range("B1").ColumnWidth = 16.109375 - $gridOffset
$ws.Range("C1").ColumnWidth = 19.44140625 - $gridOffset

# Match the author's new selection.
$ws.Range("G5").Select()
